$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first two fixtures (Man Utd v Burnley, Man Utd v Wolves) have been
# played, so remove their rows entirely; everything below shifts up.
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()
